$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the "daily" readings currently sitting in row 2 and row 4
# (row 3 is untouched) across columns D (Fecha) and M,N,O,P,Q,S,T
# (Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Unidad de comercializacion, Precio $/Kg, Kg/unidad).

$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow4 = $ws.Range($col + "4")

    $val2 = $cellRow2.Value2
    $val4 = $cellRow4.Value2

    $cellRow2.Value2 = $val4
    $cellRow4.Value2 = $val2
}
